# Auto-generated: applies cached market-price / profit recalculation values
# to the Leve profit tables on each class sheet, per scheduled-runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 348.14285
$ws.Range("I19").Value = 349.25
$ws.Range("J19").Value = 346.66666
$ws.Range("K19").Value = 349.25
$ws.Range("L19").Value = 346.66666
$ws.Range("M19").Value = -174.25
$ws.Range("N19").Value = -696.66666
$ws.Range("H33").Value = 6443.0625
$ws.Range("I33").Value = 69.7
$ws.Range("J33").Value = 17065.334
$ws.Range("K33").Value = 69.7
$ws.Range("L33").Value = 17065.334
$ws.Range("M33").Value = 159.3
$ws.Range("N33").Value = -17523.334
$ws.Range("H116").Value = 4934.273
$ws.Range("I116").Value = 4260
$ws.Range("K116").Value = 4260
$ws.Range("M116").Value = -818
$ws.Range("H132").Value = 6080.778
$ws.Range("I132").Value = 7233.3335
$ws.Range("J132").Value = 3775.6667
$ws.Range("K132").Value = 21700.0005
$ws.Range("L132").Value = 11327.0001
$ws.Range("M132").Value = -19170.0005
$ws.Range("N132").Value = -16387.0001
$ws.Range("H137").Value = 1004.5
$ws.Range("I137").Value = 801.86957
$ws.Range("K137").Value = 2405.60871
$ws.Range("M137").Value = 144.39129

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1534.2307
$ws.Range("I74").Value = 1565.7097
$ws.Range("J74").Value = 1412.25
$ws.Range("K74").Value = 1565.7097
$ws.Range("L74").Value = 1412.25
$ws.Range("M74").Value = -691.7097000000001
$ws.Range("N74").Value = -3160.25
$ws.Range("H77").Value = 1534.2307
$ws.Range("I77").Value = 1565.7097
$ws.Range("J77").Value = 1412.25
$ws.Range("K77").Value = 7828.548500000001
$ws.Range("L77").Value = 7061.25
$ws.Range("M77").Value = -3460.548500000001
$ws.Range("N77").Value = -15797.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 177.57692
$ws.Range("I80").Value = 70.166664
$ws.Range("J80").Value = 209.8
$ws.Range("K80").Value = 70.166664
$ws.Range("L80").Value = 209.8
$ws.Range("M80").Value = 927.833336
$ws.Range("N80").Value = -2205.8
$ws.Range("H83").Value = 177.57692
$ws.Range("I83").Value = 70.166664
$ws.Range("J83").Value = 209.8
$ws.Range("K83").Value = 350.83332
$ws.Range("L83").Value = 1049
$ws.Range("M83").Value = 4641.16668
$ws.Range("N83").Value = -11033

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 7325
$ws.Range("I60").Value = 4093
$ws.Range("J60").Value = 8402.333000000001
$ws.Range("K60").Value = 4093
$ws.Range("L60").Value = 8402.333000000001
$ws.Range("M60").Value = -3582
$ws.Range("N60").Value = -9424.333000000001
$ws.Range("H74").Value = 15017.429
$ws.Range("J74").Value = 17139.5
$ws.Range("L74").Value = 17139.5
$ws.Range("N74").Value = -18887.5
$ws.Range("H77").Value = 15017.429
$ws.Range("J77").Value = 17139.5
$ws.Range("L77").Value = 51418.5
$ws.Range("N77").Value = -60154.5
$ws.Range("H80").Value = 13495
$ws.Range("J80").Value = 14660
$ws.Range("L80").Value = 14660
$ws.Range("N80").Value = -16906
$ws.Range("H83").Value = 13495
$ws.Range("J83").Value = 14660
$ws.Range("L83").Value = 43980
$ws.Range("N83").Value = -55212

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1180.8462
$ws.Range("I5").Value = 804.63635
$ws.Range("J5").Value = 3250
$ws.Range("K5").Value = 2413.90905
$ws.Range("L5").Value = 9750
$ws.Range("M5").Value = -2301.90905
$ws.Range("N5").Value = -9974
$ws.Range("H57").Value = 55561376
$ws.Range("I57").Value = 500000400
$ws.Range("J57").Value = 6500
$ws.Range("K57").Value = 1500001200
$ws.Range("L57").Value = 19500
$ws.Range("M57").Value = -1500000641
$ws.Range("N57").Value = -20618
$ws.Range("H74").Value = 6984.1
$ws.Range("I74").Value = 3005.2
$ws.Range("J74").Value = 10963
$ws.Range("K74").Value = 9015.599999999999
$ws.Range("L74").Value = 32889
$ws.Range("M74").Value = -7954.599999999999
$ws.Range("N74").Value = -35011
$ws.Range("H77").Value = 6984.1
$ws.Range("I77").Value = 3005.2
$ws.Range("J77").Value = 10963
$ws.Range("K77").Value = 27046.8
$ws.Range("L77").Value = 98667
$ws.Range("M77").Value = -21742.8
$ws.Range("N77").Value = -109275
$ws.Range("H81").Value = 949.6667
$ws.Range("I81").Value = 174.5
$ws.Range("J81").Value = 2500
$ws.Range("K81").Value = 523.5
$ws.Range("L81").Value = 7500
$ws.Range("M81").Value = 599.5
$ws.Range("N81").Value = -9746
$ws.Range("H84").Value = 949.6667
$ws.Range("I84").Value = 174.5
$ws.Range("J84").Value = 2500
$ws.Range("K84").Value = 1570.5
$ws.Range("L84").Value = 22500
$ws.Range("M84").Value = 4045.5
$ws.Range("N84").Value = -33732
$ws.Range("H92").Value = 334
$ws.Range("I92").Value = 334
$ws.Range("K92").Value = 1002
$ws.Range("M92").Value = 246
$ws.Range("H135").Value = 1180.8462
$ws.Range("I135").Value = 804.63635
$ws.Range("J135").Value = 3250
$ws.Range("K135").Value = 7241.72715
$ws.Range("L135").Value = 29250
$ws.Range("M135").Value = -4706.72715
$ws.Range("N135").Value = -34320

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 2000
$ws.Range("J20").Value = 2000
$ws.Range("L20").Value = 2000
$ws.Range("N20").Value = -2490
$ws.Range("H107").Value = 413.56522
$ws.Range("I107").Value = 251.64285
$ws.Range("J107").Value = 665.44446
$ws.Range("K107").Value = 251.64285
$ws.Range("L107").Value = 665.44446
$ws.Range("M107").Value = 1668.35715
$ws.Range("N107").Value = -4505.44446
$ws.Range("H132").Value = 31526.258
$ws.Range("I132").Value = 42665.16
$ws.Range("K132").Value = 127995.48
$ws.Range("M132").Value = -125465.48

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 858857.1
$ws.Range("I2").Value = 2000
$ws.Range("J2").Value = 1001666.7
$ws.Range("K2").Value = 2000
$ws.Range("L2").Value = 1001666.7
$ws.Range("M2").Value = -1888
$ws.Range("N2").Value = -1001890.7
$ws.Range("H100").Value = 2840.3333
$ws.Range("I100").Value = 2301.5
$ws.Range("K100").Value = 2301.5
$ws.Range("M100").Value = -1760.5
$ws.Range("H132").Value = 3154.0527
$ws.Range("I132").Value = 3303.625
$ws.Range("J132").Value = 2802.1177
$ws.Range("K132").Value = 9910.875
$ws.Range("L132").Value = 8406.3531
$ws.Range("M132").Value = -7380.875
$ws.Range("N132").Value = -13466.3531

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 12249
$ws.Range("J63").Value = 12249
$ws.Range("L63").Value = 12249
$ws.Range("N63").Value = -13497
$ws.Range("H66").Value = 12249
$ws.Range("J66").Value = 12249
$ws.Range("L66").Value = 36747
$ws.Range("N66").Value = -42987
$ws.Range("H132").Value = 52639024
$ws.Range("I132").Value = 100002150
$ws.Range("J132").Value = 13329.889
$ws.Range("K132").Value = 300006450
$ws.Range("L132").Value = 39989.667
$ws.Range("M132").Value = -300003920
$ws.Range("N132").Value = -45049.667

